$p = $ppt.ActivePresentation

# 1. Update the cached "datetimeFigureOut" date text from 23/03/2023 to
#    18/03/2023 on the slide master and every slide layout's Date
#    placeholder.
$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $sh = $master.Shapes.Item($j)
    if ($sh.Name -like "*Date*") {
        $sh.TextFrame.TextRange.Text = "18/03/2023"
    }
}

$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "*Date*") {
            $sh.TextFrame.TextRange.Text = "18/03/2023"
        }
    }
}

# 2. Remove the two "Azure Cosmos DB" / "Azure Cosmos DB APIs" slides
#    (previously slides 25 and 26, sldId 282 and 283). Deleting slide 25
#    twice removes both, since the slide that used to be 26 shifts down
#    to index 25 after the first delete.
$p.Slides.Item(25).Delete()
$p.Slides.Item(25).Delete()
